$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update city names
$ws.Range("A4").Value = "Phoenix"
$ws.Range("A5").Value = "Pittsburg"

# Update row 4 (Phoenix) data
$ws.Range("B4").Value = 33.3489
$ws.Range("C4").Value = -112.4912
$ws.Range("D4").Value = 909

# Update row 5 (Pittsburg) data
$ws.Range("B5").Value = 40.4688
$ws.Range("C5").Value = -79.9812
$ws.Range("D5").Value = 886

# Update selection on the sheet view
$ws.Range("A4:D5").Select()
